$wb = $excel.ActiveWorkbook

# "zh-cn" worksheet: update Correspond Handoff / Handback Datetime values
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-20 05:21:21"
$wsZhCn.Range("E3").Value = "2016-03-20 05:21:21"
$wsZhCn.Range("H2").Value = "2016-03-20 05:22:09"
$wsZhCn.Range("H3").Value = "2016-03-20 05:22:09"

# "de-de" worksheet: update Correspond Handoff / Handback Datetime values
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-20 05:21:33"
$wsDeDe.Range("E3").Value = "2016-03-20 05:21:33"
$wsDeDe.Range("H2").Value = "2016-03-20 05:22:22"
$wsDeDe.Range("H3").Value = "2016-03-20 05:22:22"
